$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B17: ORCID ID of vocabulary creator
$ws.Range("B17").Value = "0000-0002-1665-678X"

# Update B21: modified timestamp
$ws.Range("B21").Value = "2023-06-12T17:32:46+00:00"

# New rows 27-44 data: col A = identifier, col B = label, col C = type
$newRows = @(
    @("vocab:1003", "nNOS-CreERT2 x GCaMP6f-expressing mice", "variable"),
    @("vocab:1004", "ChAT-Cre x GCaMP6f-expressing mice", "variable"),
    @("vocab:1005", "mouse study", "subject"),
    @("vocab:1006", "neurological studies", "subject"),
    @("vocab:1007", "phasic calcium transient response", "variable"),
    @("vocab:1008", "Myenteric neurons", "subject"),
    @("vocab:1009", "Murine colon", "subject"),
    @("vocab:1010", "Spontaneous colonic migrating motor complexes (CMMCs)", "subject"),
    @("vocab:1011", "Brush stimulation of the mucosa", "subject"),
    @("vocab:1012", "Electrical field stimulation (EFS)", "subject"),
    @("vocab:1013", "Elongation of the colon", "variable"),
    @("vocab:1014", "Microscopy", "subject"),
    @("vocab:1015", "Optical physiology", "subject"),
    @("vocab:1016", "Myenteric nerve plexus", "subject"),
    @("vocab:1017", "Calcium imaging", "subject"),
    @("vocab:1018", "Neuron activity", "variable"),
    @("vocab:1019", "Nikon Eclipse FN1 upright fluorescence microscope", "subject"),
    @("vocab:1020", "Photometrics Prime 95B sCMOS camera", "subject")
)

$startRow = 27
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    # Columns D through AM (4 through 39) stay blank for these rows, same as
    # the source data - the sheet's used range already spans through column
    # AM from the existing rows above, so it continues to cover them.
}
